$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (RCMP / Canada-wide scraping source)
$ws.Range("A24").Value = "CAN"
$ws.Range("B24").Value = "Canada"
$ws.Range("C24").Value = "http://www.rcmp-grc.gc.ca/detach/en/d/*/kmlloc"
$ws.Range("D24").Value = "???"
$ws.Range("E24").Value = "Needs scraping"

# Turn the source cell into a real hyperlink (adds Hyperlink style automatically)
$ws.Hyperlinks.Add($ws.Range("C24"), "http://www.rcmp-grc.gc.ca/detach/en/d/*/kmlloc")

# Widen columns C and E so the long URL / note text are fully visible
$ws.Columns.Item(3).ColumnWidth = 157.5
$ws.Columns.Item(5).ColumnWidth = 77

# Leave the selection where the user ended up after entering the new row
$ws.Range("D25").Select() | Out-Null
